$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1/E1/F1 hold values that look numeric/date ("2018-12-09", "13", "22.352")
# but must be stored as plain text, same as the rest of the row -- force
# text formatting on just those cells before assigning so Excel doesn't
# auto-convert them to a date serial / number.
$ws.Range("D1:F1").NumberFormat = "@"

$ws.Range("A1").Value = "Fountain"
$ws.Range("B1").Value = "Daurice"
$ws.Range("C1").Value = "WR"
$ws.Range("D1").Value = "2018-12-09"
$ws.Range("E1").Value = "13"
$ws.Range("F1").Value = "22.352"
$ws.Range("G1").Value = "IND"
$ws.Range("H1").Value = "@"
$ws.Range("I1").Value = "HOU"
$ws.Range("J1").Value = "W 24-21"
$ws.Range("K1").Value = ""
$ws.Range("L1").Value = 0
